$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.446.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.519.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.78%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.523.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  -3.10%  "
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.966.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.322.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.520.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0781"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("E35").Value = "  -6.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "280.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0932"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.99%  "
